# Layout.xlsx: insert a new "clue" column just before the existing column W
# (the running-index column), pushing the old W -> X. The newly inserted
# column inherits its content from the column immediately to its left (the
# last "clue" column, U), matching Excel's default insert-column behavior
# for this sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at V; this shifts the old V column to W and the old
# W (index) column to X.
$ws.Columns("V:V").Insert()

# Populate the freshly inserted column V with the same clue content as the
# column to its left (U), for every data row (1-22).
$ws.Range("U1:U22").Copy()
$ws.Range("V1:V22").PasteSpecial()

# Row 23 is a plain running index (0,1,2,...). The pre-existing row was
# missing the value 18 (S23 jumped straight from 17 to 19); fix that gap so
# the sequence is contiguous again after the column shift.
$ws.Range("S23").Value = 18
$ws.Range("T23").Value = 19
$ws.Range("U23").Value = 20
$ws.Range("V23").Value = 21

# Leave the selection on the new trailing column, one past the last used
# column (X).
$ws.Range("Y1").Select() | Out-Null
